# Add two new parameter rows to the "main" sheet:
#   mesher -> blockMesh
#   solver -> simpleFoam

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$ws.Range("A6").Value = "mesher"
$ws.Range("B6").Value = "blockMesh"

$ws.Range("A7").Value = "solver"
$ws.Range("B7").Value = "simpleFoam"

$ws.Range("B8").Select()
